$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) R10 "Resumen" cell: merge the two runs "Ubicar " + "a los entrenadores
#    y asistentes en las oficinas." into a single run with the full text.
# ---------------------------------------------------------------------------
$t9 = $d.Tables.Item(10)
$rng1 = $t9.Cell(2, 2).Range
$rng1.Find.Execute(
    "Ubicar a los entrenadores y asistentes en las oficinas.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ubicar a los entrenadores y asistentes en las oficinas.", 2)

# ---------------------------------------------------------------------------
# 2) R11 table: fill the three empty "Resumen" / "Entradas" / "Resultado"
#    cells with the new requirement text, inheriting the sz/szCs run
#    properties from the (empty) paragraph mark.
# ---------------------------------------------------------------------------
$t10 = $d.Tables.Item(11)

$rngResumen = $t10.Cell(2, 2).Range
$rngResumen.Text = "Ubicar de forma aleatoria a los jugadores en los vestuarios. En el caso del vestuario de 7x6, se debe informar que jugadores quedaron sin vestuario."
$rngResumen2 = $t10.Cell(2, 2).Range
$rngResumen2.Font.Size = 12
$rngResumen2.Font.SizeBi = 12

$rngEntradas = $t10.Cell(3, 2).Range
$rngEntradas.Text = "Team (EquipoA oEquipoB)"
$rngEntradas2 = $t10.Cell(3, 2).Range
$rngEntradas2.Font.Size = 12
$rngEntradas2.Font.SizeBi = 12

$rngResultado = $t10.Cell(4, 2).Range
$rngResultado.Text = "Se ha mostrado a los jugadores de un equipo ubicados en sus vestuarios. "
$rngResultado2 = $t10.Cell(4, 2).Range
$rngResultado2.Font.Size = 12
$rngResultado2.Font.SizeBi = 12

# ---------------------------------------------------------------------------
# 3) Repagination artifact: move <w:lastRenderedPageBreak/> from the start of
#    the R15 table's "Nombre" label run to the end of the R14 table's
#    "Resultado" label run. The cell paragraphs carry no text-level changes,
#    only the rendered-page-break marker moves, so the exact paragraph
#    identity (paraId/rsids) must be preserved.
# ---------------------------------------------------------------------------

# 3a) Add <w:lastRenderedPageBreak/> before "Resultado" text (R14 table).
$t14 = $d.Tables.Item(14)
$cellRng = $t14.Cell(4, 1).Range
$addRng = $d.Range($cellRng.Start, $cellRng.End - 1)
$addXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7B2F83E3" w14:textId="77777777" w:rsidR="00672D49" w:rsidRPr="004D3106" w:rsidRDefault="00672D49" w:rsidP="00B04461"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="004D3106"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Resultado</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$addRng.InsertXML($addXml)

# 3b) Remove <w:lastRenderedPageBreak/> from before "Nombre" text (R15 table).
$t15 = $d.Tables.Item(15)
$cellRng2 = $t15.Cell(1, 1).Range
$delRng = $d.Range($cellRng2.Start, $cellRng2.End - 1)
$delXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="114D794A" w14:textId="77777777" w:rsidR="00672D49" w:rsidRPr="004D3106" w:rsidRDefault="00672D49" w:rsidP="00B04461"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="004D3106"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Nombre</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$delRng.InsertXML($delXml)
